# Fix mojibake "Â±" (U+00C2 U+00B1, i.e. UTF-8 bytes re-interpreted/re-encoded)
# back into the correct "±" (U+00B1 PLUS-MINUS SIGN) in the metric columns
# (f1_score_weighted, training_time, test_time) of the automl results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$moji = [string][char]0x00C2 + [string][char]0x00B1   # "Â±"
$fixed = [string][char]0x00B1                          # "±"

$range = $ws.Range("B2:D17")

foreach ($cell in $range.Cells) {
    $val = $cell.Text
    if ($val -ne $null -and $val.Contains($moji)) {
        $cell.Value = $val.Replace($moji, $fixed)
    }
}
